# Recetas.xlsx - cambios en el historial del paciente
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Renombrar columna ID a "ID receta"
$ws.Range("A1").Value = "ID receta"

# Insertar 3 filas nuevas justo debajo del encabezado (filas 2,3,4)
# para las nuevas recetas, empujando las filas existentes hacia abajo.
$ws.Rows("2:4").Insert()

# Nueva fila 2: receta 7 - Carla atendida por Rosa
$ws.Range("A2").Value = "'7"
$ws.Range("B2").Value = "Carla"
$ws.Range("C2").Value = "Rosa"
$ws.Range("D2").Value = "dolocloralan"
$ws.Range("E2").Value = "2025-10-24 10:27"

# Nueva fila 3: receta 6 - José atendido por Admin
$ws.Range("A3").Value = "'6"
$ws.Range("B3").Value = "José"
$ws.Range("C3").Value = "Admin"
$ws.Range("D3").Value = "Bactrim"
$ws.Range("E3").Value = "2025-10-24 10:12"

# Nueva fila 4: receta 5 - Carla atendida por Admin
$ws.Range("A4").Value = "'5"
$ws.Range("B4").Value = "Carla"
$ws.Range("C4").Value = "Admin"
$ws.Range("D4").Value = "Migradol cada 8 horas"
$ws.Range("E4").Value = "2025-10-24 10:12"
